# Shift the "Computed Probe N start/end" window by +100 for rows 2 and 3,
# and update the downstream amplitude / ratio statistics that were
# recomputed as a result (per the diff - sweell/wind split + damping-grouped
# plotting update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("P2").Value  = 4600
$ws.Range("Q2").Value  = 4984.615384615385
$ws.Range("R2").Value  = 17.30039999999998
$ws.Range("S2").Value  = 12.40840678741196
$ws.Range("T2").Value  = 15.28740162368884

$ws.Range("AA2").Value = 4700
$ws.Range("AB2").Value = 5084.615384615385
$ws.Range("AC2").Value = 15.26759999999998
$ws.Range("AD2").Value = 11.20950937555914
$ws.Range("AE2").Value = 14.1811631901298

$ws.Range("AL2").Value = 6400
$ws.Range("AM2").Value = 6784.615384615385
$ws.Range("AN2").Value = 9.407200000000003
$ws.Range("AO2").Value = 6.698401312054257
$ws.Range("AP2").Value = 8.30968077669284

$ws.Range("AW2").Value = 6400
$ws.Range("AX2").Value = 6784.615384615385
$ws.Range("AY2").Value = 9.325400000000027
$ws.Range("AZ2").Value = 6.51790130812388
$ws.Range("BA2").Value = 8.131027086979174

$ws.Range("BK2").Value = 0.1039487496782247
$ws.Range("BR2").Value = 0.8824998265936047
$ws.Range("BS2").Value = 0.6161544709057097
$ws.Range("BT2").Value = 0.9913045326983613

# ---- Row 3 ----
$ws.Range("P3").Value  = 4600
$ws.Range("Q3").Value  = 4984.615384615385
$ws.Range("R3").Value  = 8.194999999999965
$ws.Range("S3").Value  = 6.150988042443566
$ws.Range("T3").Value  = 7.680335855956333

$ws.Range("AA3").Value = 4700
$ws.Range("AB3").Value = 5084.615384615385
$ws.Range("AC3").Value = 8.625000000000028
$ws.Range("AD3").Value = 6.572803824607096
$ws.Range("AE3").Value = 7.985192377663065

$ws.Range("AL3").Value = 6400
$ws.Range("AM3").Value = 6784.615384615385
$ws.Range("AN3").Value = 4.041800000000021
$ws.Range("AO3").Value = 2.735792981627788
$ws.Range("AP3").Value = 3.530943225084268

$ws.Range("AW3").Value = 6400
$ws.Range("AX3").Value = 6784.615384615385
$ws.Range("AY3").Value = 4.175400000000018
$ws.Range("AZ3").Value = 2.644811568780587
$ws.Range("BA3").Value = 3.371432832357906

$ws.Range("BK3").Value = 0.05872291427432548
$ws.Range("BR3").Value = 1.05247101891398
$ws.Range("BS3").Value = 0.468614492753624
$ws.Range("BT3").Value = 1.033054579642732
